$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2,8).Value = 327.5
$ws.Cells.Item(2,9).Value = 327.5
$ws.Cells.Item(2,10).Value = 0
$ws.Cells.Item(2,11).Value = 327.5
$ws.Cells.Item(2,12).Value = 0
$ws.Cells.Item(2,13).Value = -214.5
$ws.Cells.Item(2,14).ClearContents()

$ws.Cells.Item(29,8).Value = 1726.579
$ws.Cells.Item(29,9).Value = 201
$ws.Cells.Item(29,10).Value = 2271.4285
$ws.Cells.Item(29,11).Value = 603
$ws.Cells.Item(29,12).Value = 6814.2855
$ws.Cells.Item(29,13).Value = -322
$ws.Cells.Item(29,14).Value = -7376.2855

$ws.Cells.Item(38,8).Value = 1603.95
$ws.Cells.Item(38,10).Value = 2004.4839
$ws.Cells.Item(38,12).Value = 6013.4517
$ws.Cells.Item(38,14).Value = -6757.4517

$ws.Cells.Item(64,8).Value = 4716.6665
$ws.Cells.Item(64,9).Value = 5095
$ws.Cells.Item(64,10).Value = 3960
$ws.Cells.Item(64,11).Value = 5095
$ws.Cells.Item(64,12).Value = 3960
$ws.Cells.Item(64,13).Value = -4847
$ws.Cells.Item(64,14).Value = -4456

$ws.Cells.Item(67,8).Value = 4716.6665
$ws.Cells.Item(67,9).Value = 5095
$ws.Cells.Item(67,10).Value = 3960
$ws.Cells.Item(67,11).Value = 5095
$ws.Cells.Item(67,12).Value = 3960
$ws.Cells.Item(67,13).Value = -4237
$ws.Cells.Item(67,14).Value = -5676

$ws.Cells.Item(74,8).Value = 3650.375
$ws.Cells.Item(74,9).Value = 3350.75
$ws.Cells.Item(74,10).Value = 3950
$ws.Cells.Item(74,11).Value = 3350.75
$ws.Cells.Item(74,12).Value = 3950
$ws.Cells.Item(74,13).Value = -2414.75
$ws.Cells.Item(74,14).Value = -5822

$ws.Cells.Item(77,8).Value = 3650.375
$ws.Cells.Item(77,9).Value = 3350.75
$ws.Cells.Item(77,10).Value = 3950
$ws.Cells.Item(77,11).Value = 16753.75
$ws.Cells.Item(77,12).Value = 19750
$ws.Cells.Item(77,13).Value = -12073.75
$ws.Cells.Item(77,14).Value = -29110

$ws.Cells.Item(125,8).Value = 6833
$ws.Cells.Item(125,9).Value = 5249.5
$ws.Cells.Item(125,10).Value = 10000
$ws.Cells.Item(125,11).Value = 47245.5
$ws.Cells.Item(125,12).Value = 90000
$ws.Cells.Item(125,13).Value = -44785.5
$ws.Cells.Item(125,14).Value = -94920

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(63,8).Value = 2136.7354
$ws.Cells.Item(63,9).Value = 2047.3478
$ws.Cells.Item(63,10).Value = 2323.6365
$ws.Cells.Item(63,11).Value = 2047.3478
$ws.Cells.Item(63,12).Value = 2323.6365
$ws.Cells.Item(63,13).Value = -1361.3478
$ws.Cells.Item(63,14).Value = -3695.6365

$ws.Cells.Item(66,8).Value = 2136.7354
$ws.Cells.Item(66,9).Value = 2047.3478
$ws.Cells.Item(66,10).Value = 2323.6365
$ws.Cells.Item(66,11).Value = 10236.739
$ws.Cells.Item(66,12).Value = 11618.1825
$ws.Cells.Item(66,13).Value = -6804.739
$ws.Cells.Item(66,14).Value = -18482.1825

$ws.Cells.Item(74,8).Value = 2427.4546
$ws.Cells.Item(74,9).Value = 2012.75
$ws.Cells.Item(74,10).Value = 3533.3333
$ws.Cells.Item(74,11).Value = 2012.75
$ws.Cells.Item(74,12).Value = 3533.3333
$ws.Cells.Item(74,13).Value = -1138.75
$ws.Cells.Item(74,14).Value = -5281.3333

$ws.Cells.Item(77,8).Value = 2427.4546
$ws.Cells.Item(77,9).Value = 2012.75
$ws.Cells.Item(77,10).Value = 3533.3333
$ws.Cells.Item(77,11).Value = 10063.75
$ws.Cells.Item(77,12).Value = 17666.6665
$ws.Cells.Item(77,13).Value = -5695.75
$ws.Cells.Item(77,14).Value = -26402.6665

$ws.Cells.Item(88,8).Value = 2926.4546
$ws.Cells.Item(88,9).Value = 2000
$ws.Cells.Item(88,10).Value = 3019.1
$ws.Cells.Item(88,11).Value = 2000
$ws.Cells.Item(88,12).Value = 3019.1
$ws.Cells.Item(88,13).Value = -1594
$ws.Cells.Item(88,14).Value = -3831.1

$ws.Cells.Item(91,8).Value = 2926.4546
$ws.Cells.Item(91,9).Value = 2000
$ws.Cells.Item(91,10).Value = 3019.1
$ws.Cells.Item(91,11).Value = 2000
$ws.Cells.Item(91,12).Value = 3019.1
$ws.Cells.Item(91,13).Value = -596
$ws.Cells.Item(91,14).Value = -5827.1

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86,8).Value = 4916
$ws.Cells.Item(86,9).Value = 4916
$ws.Cells.Item(86,10).Value = 0
$ws.Cells.Item(86,11).Value = 4916
$ws.Cells.Item(86,12).Value = 0
$ws.Cells.Item(86,13).Value = -3793
$ws.Cells.Item(86,14).ClearContents()

$ws.Cells.Item(89,8).Value = 4916
$ws.Cells.Item(89,9).Value = 4916
$ws.Cells.Item(89,10).Value = 0
$ws.Cells.Item(89,11).Value = 24580
$ws.Cells.Item(89,12).Value = 0
$ws.Cells.Item(89,13).Value = -18964
$ws.Cells.Item(89,14).ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(62,8).Value = 33335334

$ws.Cells.Item(65,8).Value = 33335334

$ws.Cells.Item(68,8).Value = 16000
$ws.Cells.Item(68,10).Value = 16000
$ws.Cells.Item(68,12).Value = 16000
$ws.Cells.Item(68,14).Value = -17498

$ws.Cells.Item(71,8).Value = 16000
$ws.Cells.Item(71,10).Value = 16000
$ws.Cells.Item(71,12).Value = 48000
$ws.Cells.Item(71,14).Value = -55488

$ws.Cells.Item(141,8).Value = 1212986
$ws.Cells.Item(141,10).Value = 1212986
$ws.Cells.Item(141,12).Value = 1212986
$ws.Cells.Item(141,14).Value = -1223346

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(56,8).Value = 6091.9
$ws.Cells.Item(56,9).Value = 6091.9
$ws.Cells.Item(56,11).Value = 6091.9
$ws.Cells.Item(56,13).Value = -5561.9

$ws.Cells.Item(131,8).Value = 24394032
$ws.Cells.Item(131,10).Value = 4165.4053
$ws.Cells.Item(131,12).Value = 12496.2159
$ws.Cells.Item(131,14).Value = -22576.2159

$ws.Cells.Item(140,8).Value = 30175.416
$ws.Cells.Item(140,9).Value = 57703.777
$ws.Cells.Item(140,10).Value = 2647.0557
$ws.Cells.Item(140,11).Value = 173111.331
$ws.Cells.Item(140,12).Value = 7941.1671
$ws.Cells.Item(140,13).Value = -167931.331
$ws.Cells.Item(140,14).Value = -18301.1671

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80,8).Value = 3123.077
$ws.Cells.Item(80,9).Value = 1840
$ws.Cells.Item(80,11).Value = 1840
$ws.Cells.Item(80,13).Value = -842

$ws.Cells.Item(83,8).Value = 3123.077
$ws.Cells.Item(83,9).Value = 1840
$ws.Cells.Item(83,11).Value = 9200
$ws.Cells.Item(83,13).Value = -4208

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22,8).Value = 887.5
$ws.Cells.Item(22,9).Value = 549
$ws.Cells.Item(22,10).Value = 1226
$ws.Cells.Item(22,11).Value = 549
$ws.Cells.Item(22,12).Value = 1226
$ws.Cells.Item(22,13).Value = -254
$ws.Cells.Item(22,14).Value = -1816

$ws.Cells.Item(27,8).Value = 887.5
$ws.Cells.Item(27,9).Value = 549
$ws.Cells.Item(27,10).Value = 1226
$ws.Cells.Item(27,11).Value = 549
$ws.Cells.Item(27,12).Value = 1226
$ws.Cells.Item(27,13).Value = -442
$ws.Cells.Item(27,14).Value = -1440

$ws.Cells.Item(46,8).Value = 4295.4546
$ws.Cells.Item(46,9).Value = 750
$ws.Cells.Item(46,10).Value = 5625
$ws.Cells.Item(46,11).Value = 750
$ws.Cells.Item(46,12).Value = 5625
$ws.Cells.Item(46,13).Value = -562
$ws.Cells.Item(46,14).Value = -6001
